# Auto-generated edit script applying the Zodiark_Profits.xlsx diff
# Updates currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns
# (H..N) for specific rows across all 8 job sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) per the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1054.2742
$ws.Range("J17").Value = 1054.2742
$ws.Range("L17").Value = 3162.8226
$ws.Range("N17").Value = -3498.8226
# Row 62
$ws.Range("H62").Value = 9698.166999999999
$ws.Range("I62").Value = 4996.3335
$ws.Range("K62").Value = 4996.3335
$ws.Range("M62").Value = -4372.3335
# Row 65
$ws.Range("H65").Value = 9698.166999999999
$ws.Range("I65").Value = 4996.3335
$ws.Range("K65").Value = 24981.6675
$ws.Range("M65").Value = -21861.6675
# Row 76
$ws.Range("H76").Value = 7256.857
$ws.Range("J76").Value = 8224.5
$ws.Range("L76").Value = 8224.5
$ws.Range("N76").Value = -8854.5
# Row 79
$ws.Range("H79").Value = 7256.857
$ws.Range("J79").Value = 8224.5
$ws.Range("L79").Value = 8224.5
$ws.Range("N79").Value = -10408.5
# Row 95
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1764.4615
$ws.Range("I61").Value = 1598.381
$ws.Range("J61").Value = 2462
$ws.Range("K61").Value = 1598.381
$ws.Range("L61").Value = 2462
$ws.Range("M61").Value = -1386.381
$ws.Range("N61").Value = -2886
# Row 74
$ws.Range("H74").Value = 2982.2917
$ws.Range("I74").Value = 2820.2058
$ws.Range("K74").Value = 2820.2058
$ws.Range("M74").Value = -1946.2058
# Row 77
$ws.Range("H77").Value = 2982.2917
$ws.Range("I77").Value = 2820.2058
$ws.Range("K77").Value = 14101.029
$ws.Range("M77").Value = -9733.029
# Row 88
$ws.Range("H88").Value = 3513.111
$ws.Range("J88").Value = 3598.1428
$ws.Range("L88").Value = 3598.1428
$ws.Range("N88").Value = -4410.1428
# Row 91
$ws.Range("H91").Value = 3513.111
$ws.Range("J91").Value = 3598.1428
$ws.Range("L91").Value = 3598.1428
$ws.Range("N91").Value = -6406.1428
# Row 119
$ws.Range("H119").Value = 49950
$ws.Range("J119").Value = 49950
$ws.Range("L119").Value = 49950
$ws.Range("N119").Value = -59626
# Row 136
$ws.Range("H136").Value = 1764.4615
$ws.Range("I136").Value = 1598.381
$ws.Range("J136").Value = 2462
$ws.Range("K136").Value = 4795.143
$ws.Range("L136").Value = 7386
$ws.Range("M136").Value = -2245.143
$ws.Range("N136").Value = -12486

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3357.2727
$ws.Range("I86").Value = 2625.1667
$ws.Range("J86").Value = 4235.8
$ws.Range("K86").Value = 2625.1667
$ws.Range("L86").Value = 4235.8
$ws.Range("M86").Value = -1502.1667
$ws.Range("N86").Value = -6481.8
# Row 89
$ws.Range("H89").Value = 3357.2727
$ws.Range("I89").Value = 2625.1667
$ws.Range("J89").Value = 4235.8
$ws.Range("K89").Value = 13125.8335
$ws.Range("L89").Value = 21179
$ws.Range("M89").Value = -7509.833500000001
$ws.Range("N89").Value = -32411
# Row 105
$ws.Range("H105").Value = 3080.9375
$ws.Range("I105").Value = 3027.125
$ws.Range("J105").Value = 3134.75
$ws.Range("K105").Value = 3027.125
$ws.Range("L105").Value = 3134.75
$ws.Range("M105").Value = -1280.125
$ws.Range("N105").Value = -6628.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2566.3
$ws.Range("I31").Value = 5705
$ws.Range("K31").Value = 5705
$ws.Range("M31").Value = -5410
# Row 34
$ws.Range("H34").Value = 2566.3
$ws.Range("I34").Value = 5705
$ws.Range("K34").Value = 5705
$ws.Range("M34").Value = -5503
# Row 105
$ws.Range("H105").Value = 23313.75
$ws.Range("I105").Value = 59403.332
$ws.Range("K105").Value = 59403.332
$ws.Range("M105").Value = -57656.332

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 802.61536
$ws.Range("I7").Value = 774.9
$ws.Range("J7").Value = 895
$ws.Range("K7").Value = 2324.7
$ws.Range("L7").Value = 2685
$ws.Range("M7").Value = -2212.7
$ws.Range("N7").Value = -2909
# Row 35
$ws.Range("H35").Value = 999
$ws.Range("J35").Value = 999
$ws.Range("L35").Value = 2997
$ws.Range("N35").Value = -3573
# Row 36
$ws.Range("H36").Value = 336
$ws.Range("I36").Value = 336
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1008
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -839
$ws.Range("N36").ClearContents()
# Row 38
$ws.Range("H38").Value = 155.83333
$ws.Range("J38").Value = 176.4
$ws.Range("L38").Value = 529.2
$ws.Range("N38").Value = -1223.2
# Row 41
$ws.Range("H41").Value = 124
$ws.Range("I41").Value = 98
$ws.Range("K41").Value = 294
$ws.Range("M41").Value = 44
# Row 42
$ws.Range("H42").Value = 2500
$ws.Range("J42").Value = 2500
$ws.Range("L42").Value = 7500
$ws.Range("N42").Value = -8568
# Row 44
$ws.Range("H44").Value = 473.75
$ws.Range("I44").Value = 466
$ws.Range("J44").Value = 497
$ws.Range("K44").Value = 1398
$ws.Range("L44").Value = 1491
$ws.Range("M44").Value = -1000
$ws.Range("N44").Value = -2287
# Row 49
$ws.Range("H49").Value = 3891.5
$ws.Range("I49").Value = 2503
$ws.Range("K49").Value = 7509
$ws.Range("M49").Value = -7353
# Row 55
$ws.Range("H55").Value = 4289.9
$ws.Range("J55").Value = 10333.333
$ws.Range("L55").Value = 30999.999
$ws.Range("N55").Value = -31353.999
# Row 69
$ws.Range("H69").Value = 8029.353
$ws.Range("I69").Value = 8281.25
$ws.Range("K69").Value = 24843.75
$ws.Range("M69").Value = -24032.75
# Row 72
$ws.Range("H72").Value = 8029.353
$ws.Range("I72").Value = 8281.25
$ws.Range("K72").Value = 74531.25
$ws.Range("M72").Value = -70475.25
# Row 74
$ws.Range("H74").Value = 11666.5
$ws.Range("I74").Value = 7499.75
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 22499.25
$ws.Range("L74").Value = 60000
$ws.Range("M74").Value = -21438.25
$ws.Range("N74").Value = -62122
# Row 77
$ws.Range("H77").Value = 11666.5
$ws.Range("I77").Value = 7499.75
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 67497.75
$ws.Range("L77").Value = 180000
$ws.Range("M77").Value = -62193.75
$ws.Range("N77").Value = -190608

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
# Row 47
$ws.Range("H47").Value = 199993.33
$ws.Range("J47").Value = 199993.33
$ws.Range("L47").Value = 199993.33
$ws.Range("N47").Value = -201129.33

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 4373.25
$ws.Range("I68").Value = 1445.1666
$ws.Range("K68").Value = 1445.1666
$ws.Range("M68").Value = -696.1666
# Row 71
$ws.Range("H71").Value = 4373.25
$ws.Range("I71").Value = 1445.1666
$ws.Range("K71").Value = 7225.833000000001
$ws.Range("M71").Value = -3481.833000000001
# Row 93
$ws.Range("H93").Value = 3222.8462
$ws.Range("I93").Value = 1478.2
$ws.Range("J93").Value = 9038.333000000001
$ws.Range("K93").Value = 1478.2
$ws.Range("L93").Value = 9038.333000000001
$ws.Range("M93").Value = -230.2
$ws.Range("N93").Value = -11534.333
# Row 101
$ws.Range("H101").Value = 85232.2
$ws.Range("J101").Value = 85232.2
$ws.Range("L101").Value = 85232.2
$ws.Range("N101").Value = -91722.2
# Row 122
$ws.Range("H122").Value = 5534.0527
$ws.Range("I122").Value = 4264.7
$ws.Range("K122").Value = 12794.1
$ws.Range("M122").Value = -10344.1
# Row 136
$ws.Range("H136").Value = 4126.25
$ws.Range("I136").Value = 3687.75
$ws.Range("K136").Value = 11063.25
$ws.Range("M136").Value = -8513.25

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3623.9375
$ws.Range("J81").Value = 4452.7
$ws.Range("L81").Value = 8905.4
$ws.Range("N81").Value = -11027.4
# Row 84
$ws.Range("H84").Value = 3623.9375
$ws.Range("J84").Value = 4452.7
$ws.Range("L84").Value = 44527
$ws.Range("N84").Value = -55135
